$d = $word.ActiveDocument
$ps = $d.PageSetup
$ps.TopMargin = $ps.TopMargin
